$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 96, shifting rows 96:149 down to 97:150.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new weekly reading.
$ws.Cells.Item(96, 1).Value2 = 4
$ws.Cells.Item(96, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(96, 3).Value2 = "Los Lagos"
$ws.Cells.Item(96, 4).Value2 = 44460
$ws.Cells.Item(96, 5).Value2 = 10
$ws.Cells.Item(96, 6).Value2 = 100112044
$ws.Cells.Item(96, 7).Value2 = "Perejil"
$ws.Cells.Item(96, 8).Value2 = "Sin especificar"
$ws.Cells.Item(96, 9).Value2 = "Primera"
$ws.Cells.Item(96, 10).Value2 = 90
$ws.Cells.Item(96, 11).Value2 = 6500
$ws.Cells.Item(96, 12).Value2 = 6500
$ws.Cells.Item(96, 13).Value2 = 6500
$ws.Cells.Item(96, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(96, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(96, 16).Value2 = 2167
$ws.Cells.Item(96, 17).Value2 = 3
$ws.Cells.Item(96, 18).Value2 = "Hortaliza"

# Keep the date formatting style consistent with the other date cells (numFmt 165).
$ws.Cells.Item(96, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
